$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 318; existing rows 318-420 shift down to 319-421.
$ws.Rows(318).Insert()

# Populate the newly inserted row 318 with the new weekly data point.
$ws.Range("A318").Value = 5
$ws.Range("B318").Value = "Macroferia Regional de Talca"
$ws.Range("C318").Value = "Maule"
$ws.Range("D318").Value = 44627
$ws.Range("E318").Value = 7
$ws.Range("F318").Value = 100114001
$ws.Range("G318").Value = "Papa"
$ws.Range("H318").Value = "Patagonia"
$ws.Range("I318").Value = "1a (cosecha)"
$ws.Range("J318").Value = 1300
$ws.Range("K318").Value = 6000
$ws.Range("L318").Value = 6000
$ws.Range("M318").Value = 6000
$ws.Range("N318").Value = "$/saco 25 kilos"
$ws.Range("O318").Value = "Región de Los Lagos"
$ws.Range("P318").Value = 240
$ws.Range("Q318").Value = 25
$ws.Range("R318").Value = "Hortaliza"
